$wb = $excel.ActiveWorkbook

# --- Insert the new "placeOrderGuest" sheet between "placeOrder" and
# --- "placeMultipleItems" (mirrors Insert Worksheet > rename, done before
# --- touching the other two sheets so the active-tab bookkeeping matches).
$anchor = $wb.Worksheets.Item("placeMultipleItems")
$guest = $wb.Worksheets.Add($anchor)
$guest.Name = "placeOrderGuest"

# Header row
$guest.Range("A1").Value = "searchKey"
$guest.Range("B1").Value = "Product Name"
$guest.Range("C1").Value = "Quantity"
$guest.Range("D1").Value = "BillingAndDeliveryAddressSame"

# Row 2 - ipod / iPod Classic / 1 / true (quantity & flag stored as text,
# same convention the existing sheets use for their Quantity column)
$guest.Range("A2").Value = "ipod"
$guest.Range("B2").Value = "iPod Classic"
$guest.Range("C2").Value = "'1"
$guest.Range("D2").Value = "'true"

# Row 3 - sony / Sony VAIO / 2 / false
$guest.Range("A3").Value = "sony"
$guest.Range("B3").Value = "Sony VAIO"
$guest.Range("C3").Value = "'2"
$guest.Range("D3").Value = "'false"

# --- Selections left behind on each sheet, matching the recorded state ---
$placeOrder = $wb.Worksheets.Item("placeOrder")
$placeOrder.Range("A3:C3").Select() | Out-Null

$placeMultipleItems = $wb.Worksheets.Item("placeMultipleItems")
$placeMultipleItems.Range("A4:C4").Select() | Out-Null

# Leave the new guest-order sheet active/selected last, as the recorded
# workbook has it as the active tab.
$guest.Activate() | Out-Null
$guest.Range("C3").Select() | Out-Null
